# "Solucion a taller de modulo R"
#
# The author finished the exercise by filling in the missing
# "factor_frec" value for the last row (cos/sen table) and then
# formatted the factor_frec column (C) as a plain integer number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing factor_frec value on the last data row.
$ws.Range("C7").Value = 1

# Apply an integer ("0") number format to the whole factor_frec column,
# like the rest of the numeric values in that column.
$ws.Columns("C").NumberFormat = "0"

# Leave the cursor where the author left it when finishing up.
$ws.Range("D14").Select()

Write-Output "Applied solution edits to ejemplo2 worksheet"
